$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 58256
$ws.Range("B3").Value = 91766
$ws.Range("B4").Value = 57830
$ws.Range("B5").Value = 58043
